# Rearrange module D to calculate emissions from activity and ef.
# Add a "type" column (comb / NC) to the Sectors sheet so that the
# activity and emission-factor databases line up for both combustion
# and non-combustion (process) emissions.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sectors")

# --- Header: D1 = "type", formatted like the other header cells (C1) ---
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D1").Value = "type"

# --- Body rows: D2:D29 -> "comb" (sector units = kt, combustion sectors) ---
#     D30:D59 -> "NC"   (sector units = B2005USD, non-combustion sectors) ---
# Use the existing D2 blank-but-formatted cell as the format source for the
# body rows (D37 currently has the wrong style, so normalise it too).
$ws.Range("D2").Copy()
$ws.Range("D2:D59").PasteSpecial(-4122)

for ($row = 2; $row -le 29; $row++) {
    $ws.Cells.Item($row, 4).Value = "comb"
}

for ($row = 30; $row -le 59; $row++) {
    $ws.Cells.Item($row, 4).Value = "NC"
}

# --- Selection state: sheet is scrolled back to the top and D5 is selected ---
$ws.Activate()
$ws.Range("D5").Select()
